# Project Sample Project is saved. Author: admin. Type: SAVE.
#
# Change: cell B11 on the "Rules" sheet used to hold the label "R40"
# (the 4th rule-row marker). It is retyped as the text value "1",
# stored as a shared string (not a number) so the cell keeps displaying
# and round-tripping as text, while its existing cell format/style is
# left completely untouched.
#
# Forcing Excel to store a numeric-looking string ("1") as text
# normally requires flipping the cell's number format to "Text" for the
# moment of entry - which would otherwise stick around afterwards and
# mint a brand new cell style for B11. To avoid mutating B11's original
# style, do the "type as text" dance on a scratch cell far outside the
# sheet's used range, then bring only the resulting text value back
# onto B11, restoring B11's original formatting afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target  = $ws.Range("B11")
$scratch = $ws.Range("Z100")

# Remember B11's current formatting on the scratch cell so it can be
# pasted back after the value swap below.
$target.Copy($scratch)

# Type "1" into B11 as text: switch to a text number format so the
# numeric-looking literal is stored as a shared string (t="s") rather
# than a number, then restore the cell's real formatting on top of
# that (pure values stay, formatting reverts to the original).
$target.NumberFormat = "@"
$target.Value = "1"

$scratch.Copy()
$target.PasteSpecial(-4122)  # xlPasteFormats

# Clean up the scratch cell so nothing else in the sheet changes.
$scratch.Clear()
